{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2024-12-10 Tuesday\", \"2024-12-11 Wednesday\"],\n  [\"700\u00d77=\", \"858\u00d77=\"],\n  [\"583\u00d73=\", \"472\u00d78=\"],\n  [\"528\u00d72=\", \"449\u00d76=\"],\n  [\"997\u00d75=\", \"116\u00d73=\"],\n  [\"122\u00d78=\", \"623\u00d75=\"],\n  [\"708\u00d73=\", \"681\u00d75=\"],\n  [\"612\u00d72=\", \"529\u00d72=\"],\n  [\"991\u00d79=\", \"580\u00d74=\"],\n  [\"489\u00d79=\", \"475\u00d76=\"],\n  [\"987\u00d78=\", \"777\u00d78=\"],\n  [\"658\u00d72=\", \"871\u00d72=\"],\n  [\"602\u00d76=\", \"825\u00d75=\"],\n  [\"765\u00d76=\", \"871\u00d79=\"],\n  [\"260\u00d79=\", \"217\u00d76=\"],\n  [\"417\u00d78=\", \"329\u00d76=\"],\n  [\"432\u00d73=\", \"317\u00d74=\"],\n  [\"286\u00d73=\", \"734\u00d79=\"],\n  [\"114\u00d76=\", \"982\u00d78=\"],\n  [\"340\u00d77=\", \"828\u00d72=\"],\n  [\"381\u00d75=\", \"301\u00d77=\"],\n  [\"130\u00d72=\", \"297\u00d74=\"],\n  [\"192\u00d78=\", \"403\u00d73=\"],\n  [\"594\u00d73=\", \"679\u00d77=\"],\n  [\"452\u00d79=\", \"348\u00d72=\"],\n  [\"411\u00d72=\", \"639\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, 'Replace');\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each pair is [old text to find, new text to put in its place]. A leading\n# \",\" on every array literal keeps PowerShell from flattening the list of\n# pairs into one long flat list of strings.\n$replacements = @(\n    ,@(\"2024-12-10 Tuesday\", \"2024-12-11 Wednesday\")\n    ,@(\"700\u00d77=\", \"858\u00d77=\")\n    ,@(\"583\u00d73=\", \"472\u00d78=\")\n    ,@(\"528\u00d72=\", \"449\u00d76=\")\n    ,@(\"997\u00d75=\", \"116\u00d73=\")\n    ,@(\"122\u00d78=\", \"623\u00d75=\")\n    ,@(\"708\u00d73=\", \"681\u00d75=\")\n    ,@(\"612\u00d72=\", \"529\u00d72=\")\n    ,@(\"991\u00d79=\", \"580\u00d74=\")\n    ,@(\"489\u00d79=\", \"475\u00d76=\")\n    ,@(\"987\u00d78=\", \"777\u00d78=\")\n    ,@(\"658\u00d72=\", \"871\u00d72=\")\n    ,@(\"602\u00d76=\", \"825\u00d75=\")\n    ,@(\"765\u00d76=\", \"871\u00d79=\")\n    ,@(\"260\u00d79=\", \"217\u00d76=\")\n    ,@(\"417\u00d78=\", \"329\u00d76=\")\n    ,@(\"432\u00d73=\", \"317\u00d74=\")\n    ,@(\"286\u00d73=\", \"734\u00d79=\")\n    ,@(\"114\u00d76=\", \"982\u00d78=\")\n    ,@(\"340\u00d77=\", \"828\u00d72=\")\n    ,@(\"381\u00d75=\", \"301\u00d77=\")\n    ,@(\"130\u00d72=\", \"297\u00d74=\")\n    ,@(\"192\u00d78=\", \"403\u00d73=\")\n    ,@(\"594\u00d73=\", \"679\u00d77=\")\n    ,@(\"452\u00d79=\", \"348\u00d72=\")\n    ,@(\"411\u00d72=\", \"639\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace) \u2014 wdFindContinue(1), wdReplaceAll(2).\n    $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
